$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value2 = '25.790.68'
$ws.Range('E2').Value2 = '  -0.17%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value2 = '1.639.27'
$ws.Range('E3').Value2 = '  +0.36%  '
$ws.Range('E4').Value2 = '  -0.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value2 = '215.77'
$ws.Range('E5').Value2 = '  +0.43%  '
$ws.Range('E7').Value2 = '  -0.07%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value2 = '0.257'
$ws.Range('E8').Value2 = '  -0.07%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value2 = '0.0635'
$ws.Range('E9').Value2 = '  -1.01%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value2 = '19.66'
$ws.Range('E10').Value2 = '  -1.03%  '
$ws.Range('E11').Value2 = '  +1.38%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value2 = '1.865.42'
$ws.Range('E13').Value2 = '  +0.34%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value2 = '1.639.89'
$ws.Range('E14').Value2 = '  +0.31%  '
$ws.Range('E15').Value2 = '  +0.77%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value2 = '0.0₃0765'
$ws.Range('E16').Value2 = '  -0.26%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value2 = '63.16'
$ws.Range('E17').Value2 = '  +0.34%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value2 = '25.842.61'
$ws.Range('E18').Value2 = '  -0.03%  '
$ws.Range('E19').Value2 = '  -0.07%  '
$ws.Range('E20').Value2 = '  +2.07%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value2 = '192.36'
$ws.Range('E21').Value2 = '  -0.63%  '
$ws.Range('E22').Value2 = '  +0.62%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value2 = '6.32'
$ws.Range('E23').Value2 = '  +1.48%  '
$ws.Range('E24').Value2 = '  +7.07%  '
$ws.Range('E25').Value2 = '  -0.01%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value2 = '142.02'
$ws.Range('E26').Value2 = '  +2.08%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value2 = '0.124'
$ws.Range('E27').Value2 = '  +1.38%  '
$ws.Range('E28').Value2 = '  +1.59%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value2 = '15.50'
$ws.Range('E30').Value2 = '  +0.16%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value2 = '0.0492'
$ws.Range('E31').Value2 = '  -0.31%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value2 = '3.32'
$ws.Range('E32').Value2 = '  +0.84%  '
$ws.Range('E33').Value2 = '  -0.34%  '
$ws.Range('E34').Value2 = '  +0.30%  '
$ws.Range('E35').Value2 = '  -0.08%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value2 = '0.906'
$ws.Range('E36').Value2 = '  +0.56%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value2 = '1.134.36'
$ws.Range('E37').Value2 = '  +1.07%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value2 = '2.54'
$ws.Range('E38').Value2 = '  -1.66%  '
$ws.Range('E39').Value2 = '  -0.80%  '
$ws.Range('E40').Value2 = '  +0.06%  '
$ws.Range('E41').Value2 = '  +0.15%  '
$ws.Range('E42').Value2 = '  +0.81%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value2 = '100.69'
$ws.Range('E43').Value2 = '  +1.06%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value2 = '0.807'
$ws.Range('E44').Value2 = '  +0.63%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value2 = '1.774.85'
$ws.Range('E45').Value2 = '  -0.02%  '
$ws.Range('E46').Value2 = '  +2.63%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value2 = '55.31'
$ws.Range('E47').Value2 = '  -0.31%  '
$ws.Range('E48').Value2 = '  -1.24%  '
$ws.Range('E49').Value2 = '  -0.23%  '
$ws.Range('E50').Value2 = '  +4.16%  '
$ws.Range('B51').Value2 = 'Algorand'
$ws.Range('C51').Value2 = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value2 = '0.0957'
$ws.Range('E51').Value2 = '  +2.13%  '
